$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3491903333333333
$ws.Range("H2").Value = 1.047571
$ws.Range("I2").Value = 0.008130334326258625
$ws.Range("J2").Value = 0.008130334326258625
$ws.Range("M2").Value = 139.728498
$ws.Range("N2").Value = 419.185494
$ws.Range("O2").Value = 0.9065295391216045
$ws.Range("P2").Value = 0.9065295391216045
$ws.Range("Q2").Value = 48.791840792786
$ws.Range("R2").Value = 439.126567135074
$ws.Range("S2").Value = 0.007370388229687793
$ws.Range("T2").Value = 0.007370388229687793

$ws.Range("G3").Value = 0.3491903333333333
$ws.Range("H3").Value = 1.047571
$ws.Range("I3").Value = 0.008130334326258625
$ws.Range("J3").Value = 0.008130334326258625
$ws.Range("O3").Value = 0.005362677585431591
$ws.Range("P3").Value = 0.005362677585431591
$ws.Range("Q3").Value = 0.288633629329889
$ws.Range("R3").Value = 2.597702663969001
$ws.Range("S3").Value = 0.00004360036165349218
$ws.Range("T3").Value = 0.00004360036165349218

$ws.Range("G4").Value = 0.3491903333333333
$ws.Range("H4").Value = 1.047571
$ws.Range("I4").Value = 0.008130334326258625
$ws.Range("J4").Value = 0.008130334326258625
$ws.Range("O4").Value = 0.0881077832929639
$ws.Range("P4").Value = 0.0881077832929639
$ws.Range("Q4").Value = 4.742196199366111
$ws.Range("R4").Value = 42.679765794295
$ws.Range("S4").Value = 0.0007163457349173406
$ws.Range("T4").Value = 0.0007163457349173406

$ws.Range("I5").Value = 0.801301577139928
$ws.Range("J5").Value = 0.8013015771399279
$ws.Range("M5").Value = 139.728498
$ws.Range("N5").Value = 419.185494
$ws.Range("O5").Value = 0.9065295391216045
$ws.Range("P5").Value = 0.9065295391216045
$ws.Range("Q5").Value = 4808.778754958179
$ws.Range("R5").Value = 43279.00879462361
$ws.Range("S5").Value = 0.7264035494220737
$ws.Range("T5").Value = 0.7264035494220736

$ws.Range("I6").Value = 0.801301577139928
$ws.Range("J6").Value = 0.8013015771399279
$ws.Range("O6").Value = 0.005362677585431591
$ws.Range("P6").Value = 0.005362677585431591
$ws.Range("S6").Value = 0.004297122006899275
$ws.Range("T6").Value = 0.004297122006899274

$ws.Range("I7").Value = 0.801301577139928
$ws.Range("J7").Value = 0.8013015771399279
$ws.Range("O7").Value = 0.0881077832929639
$ws.Range("P7").Value = 0.0881077832929639
$ws.Range("S7").Value = 0.07060090571095497
$ws.Range("T7").Value = 0.07060090571095497

$ws.Range("I8").Value = 0.1905680885338134
$ws.Range("J8").Value = 0.1905680885338134
$ws.Range("M8").Value = 139.728498
$ws.Range("N8").Value = 419.185494
$ws.Range("O8").Value = 0.9065295391216045
$ws.Range("P8").Value = 0.9065295391216045
$ws.Range("Q8").Value = 1143.639051336054
$ws.Range("R8").Value = 10292.75146202449
$ws.Range("S8").Value = 0.172755601469843
$ws.Range("T8").Value = 0.172755601469843

$ws.Range("I9").Value = 0.1905680885338134
$ws.Range("J9").Value = 0.1905680885338134
$ws.Range("O9").Value = 0.005362677585431591
$ws.Range("P9").Value = 0.005362677585431591
$ws.Range("S9").Value = 0.001021955216878824
$ws.Range("T9").Value = 0.001021955216878824

$ws.Range("I10").Value = 0.1905680885338134
$ws.Range("J10").Value = 0.1905680885338134
$ws.Range("O10").Value = 0.0881077832929639
$ws.Range("P10").Value = 0.0881077832929639
$ws.Range("S10").Value = 0.01679053184709159
$ws.Range("T10").Value = 0.01679053184709159

